# "updated roll and grades"
# Fill in the missing Assignment marks for a few students and
# leave the grid positioned/zoomed where the user last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lim, Brian (row 14) also completed Assignment 2
$ws.Range("E14").Value = 1

# Trejo, Vincent A. (row 26) also completed Assignment 1
$ws.Range("D26").Value = 1

# Vig, Mrinal K. (row 27) also completed Assignment 2
$ws.Range("E27").Value = 1

# Update the sheet view: selected cell and zoom level
[void]$ws.Range("C32").Select()
$excel.ActiveWindow.Zoom = 85
